$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Replace "LM358 Op-amp" with "MCP6002 Op-amp" in the POWER section
#    (row 16) and update its JLCPCB code + unit cost.
# ------------------------------------------------------------------
$ws.Range("B16").Value = "MCP6002 Op-amp"
$ws.Range("C16").Value = "C7377"
$ws.Range("E16").Value = 0.4035

# ------------------------------------------------------------------
# 2) Update several component unit costs in the INTERFACING section
#    (their name/part-code stay the same, only the price changes).
# ------------------------------------------------------------------
$ws.Range("E24").Value = 0.0296          # Green LED
$ws.Range("E25").Value = 0.001           # 100nF Capacitor
$ws.Range("E26").Value = 0.0055          # 1.0uF Capacitor
$ws.Range("E30").Value = 0.0013          # 1.2kOhm Resistor
$ws.Range("E31").Value = 0.0026          # 4.7kOhm Resistor

# ------------------------------------------------------------------
# 3) Insert a brand-new row right after row 31 (before the
#    INTERFACING "Total" row) for the 10kOhm Resistor, pushing the
#    Total row, the blank spacer and the whole SENSING block down by
#    one row - exactly what Excel's native "Insert Row" does.
# ------------------------------------------------------------------
$ws.Rows("32").Insert()

$ws.Range("B32").Value = "10kΩ Resistor"
$ws.Range("C32").Value = "C25804"
$ws.Range("D32").Value = 2
$ws.Range("E32").Value = 0.0012
$ws.Range("F32").Formula = "=(D32*E32)"

# Match the formatting of the row above (D column vertical-top style)
$ws.Range("D31").Copy()
$ws.Range("D32").PasteSpecial(-4122)

# Extend the INTERFACING "Total" row's SUM range to pick up the new row.
$ws.Range("F33").Formula = "=SUM(F23:F32)"
